# Weekly update: insert 3 new rows (Banquete/Primera/Segunda, "Sin especificar")
# for the newest reporting date (2021-10-25, serial 44494) above the existing
# "Mercado Mayorista Lo Valledor de Santiago - Espárragos" data, pushing the
# prior rows 35-65 down to 38-68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 35; everything from row 35 down shifts to row 38.
$ws.Rows("35:37").Insert()

# Common, repeated values for the new rows.
$mercadoId = 6
$mercado   = "Mercado Mayorista Lo Valledor de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$catId     = 300000000
$categoria = "Espárragos"
$variedad  = "Sin especificar"
$fecha     = 44494
$unidad    = "`$/kilo"
$kgUnid    = 1
$clasif    = "Hortaliza"

# row, calidad, volumen, precioMin, precioMax, precioProm, origen
$rows = @(
    @(35, "Banquete", 970, 1300, 1400, 1367, "Provincia de Linares"),
    @(36, "Primera",  880, 1100, 1200, 1151, "Provincia de Linares"),
    @(37, "Segunda",  510,  900, 1000,  949, "Provincia de Linares")
)

foreach ($r in $rows) {
    $row        = $r[0]
    $calidad    = $r[1]
    $volumen    = $r[2]
    $precioMin  = $r[3]
    $precioMax  = $r[4]
    $precioProm = $r[5]
    $origen     = $r[6]

    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $catId
    $ws.Cells.Item($row, 7).Value2  = $categoria
    $ws.Cells.Item($row, 8).Value2  = $variedad
    $ws.Cells.Item($row, 9).Value2  = $calidad
    $ws.Cells.Item($row, 10).Value2 = $volumen
    $ws.Cells.Item($row, 11).Value2 = $precioMin
    $ws.Cells.Item($row, 12).Value2 = $precioMax
    $ws.Cells.Item($row, 13).Value2 = $precioProm
    $ws.Cells.Item($row, 14).Value2 = $unidad
    $ws.Cells.Item($row, 15).Value2 = $origen
    $ws.Cells.Item($row, 16).Value2 = $precioProm
    $ws.Cells.Item($row, 17).Value2 = $kgUnid
    $ws.Cells.Item($row, 18).Value2 = $clasif
}
